# Rename the three inline "logo" pictures that live in the document's
# headers/footers:
#   - the Pearson logo picture embedded in both footers: image1.png -> image2.png
#   - the BTEC logo picture embedded in the (first-page) header: image2.jpg -> image1.jpg
#
# Word does not let you address these drawings from the body's
# InlineShapes collection (they live in header/footer stories), so we
# reach them via Sections(1).Headers/Footers(...).Range.InlineShapes.
# Renaming an InlineShape has to go through the Selection object (a
# freshly fetched InlineShape reference rejects a direct .Name write),
# so we select the picture's own Range first and then rename it through
# $word.Selection.InlineShapes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($rng, $newName) {
    $shape = $rng.InlineShapes.Item(1)
    $shape.Range.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Primary footer: Pearson logo, name image1.png -> image2.png
Rename-InlineLogo $sec.Footers.Item(1).Range "image2.png"

# First-page footer: Pearson logo, name image1.png -> image2.png
Rename-InlineLogo $sec.Footers.Item(2).Range "image2.png"

# First-page header: BTEC logo, name image2.jpg -> image1.jpg
Rename-InlineLogo $sec.Headers.Item(2).Range "image1.jpg"
